$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-23 Saturday" "2024-11-24 Sunday"

Replace-Text "218÷7=31, 1" "152÷2=76, 0"
Replace-Text "540÷3=180, 0" "303÷2=151, 1"
Replace-Text "165÷4=41, 1" "485÷8=60, 5"
Replace-Text "249÷9=27, 6" "752÷8=94, 0"
Replace-Text "804÷9=89, 3" "429÷9=47, 6"

Replace-Text "985÷8=123, 1" "170÷7=24, 2"
Replace-Text "344÷9=38, 2" "388÷7=55, 3"
Replace-Text "221÷2=110, 1" "291÷2=145, 1"
Replace-Text "134÷2=67, 0" "869÷4=217, 1"
Replace-Text "192÷9=21, 3" "494÷2=247, 0"

Replace-Text "456÷3=152, 0" "586÷6=97, 4"
Replace-Text "910÷5=182, 0" "776÷6=129, 2"
Replace-Text "543÷8=67, 7" "612÷3=204, 0"
Replace-Text "286÷5=57, 1" "711÷9=79, 0"
Replace-Text "576÷7=82, 2" "318÷4=79, 2"

Replace-Text "655÷4=163, 3" "155÷6=25, 5"
Replace-Text "480÷7=68, 4" "135÷8=16, 7"
Replace-Text "939÷7=134, 1" "848÷8=106, 0"
Replace-Text "863÷7=123, 2" "762÷8=95, 2"
Replace-Text "270÷2=135, 0" "145÷2=72, 1"

Replace-Text "172÷9=19, 1" "281÷7=40, 1"
Replace-Text "913÷5=182, 3" "757÷5=151, 2"
Replace-Text "909÷7=129, 6" "512÷5=102, 2"
Replace-Text "675÷7=96, 3" "995÷5=199, 0"
Replace-Text "128÷9=14, 2" "407÷9=45, 2"

$d.Save()
